$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the last ("Then / the result should be") row,
# pushing it from row 5 down to row 6.
$ws.Rows("5").Insert()

# Keep selection in sync with the shifted formula cell.
$ws.Range("C6").Select()
